$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.243.69'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '1.674.90'
$ws.Range('E3').Value = '  +2.85%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '219.62'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.42%  '
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '29.70'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.84%  '
$ws.Range('E9').Value = '  +2.18%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0619'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.96%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0906'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.06%  '
$ws.Range('D12').Value = '1.914.18'
$ws.Range('E12').Value = '  +2.78%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.75'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +19.16%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.620'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +9.07%  '
$ws.Range('D15').Value = '1.678.07'
$ws.Range('E15').Value = '  +2.99%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.01'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.75%  '
$ws.Range('D17').Value = '30.244.58'
$ws.Range('E17').Value = '  +1.06%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '65.72'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '247.47'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('D20').Value = '0.0₃0719'
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.33'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +4.91%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.05'
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.22'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +5.26%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '158.89'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '15.88'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.15%  '
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.77'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.51%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0501'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.07%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.48'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +3.81%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.13'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.30'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.43%  '
$ws.Range('D34').Value = '1.479.61'
$ws.Range('E34').Value = '  +3.55%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.74'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +5.81%  '
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0179'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +4.72%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.590'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +6.21%  '
$ws.Range('B39').Value = 'Aave'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '79.61'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +15.26%  '
$ws.Range('E40').Value = '  -7.00%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.30'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.859'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.18%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.01'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +2.54%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0505'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('E45').Value = '  -2.91%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.999'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '52.17'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -4.84%  '
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').Value = '1.807.96'
$ws.Range('E49').Value = '  +2.11%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '95.54'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +6.49%  '
$ws.Range('E51').Value = '  +10.18%  '
